function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $c = $Sheet.Range($Addr)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 'D2' '42.739.75'
Set-TextValue $ws 'E2' '  -1.63%  '
Set-TextValue $ws 'D3' '2.306.62'
Set-TextValue $ws 'E3' '  -0.22%  '
Set-TextValue $ws 'E4' '  +0.06%  '
Set-TextValue $ws 'D5' '304.14'
Set-TextValue $ws 'E5' '  -2.18%  '
Set-TextValue $ws 'D6' '99.34'
Set-TextValue $ws 'E6' '  -4.60%  '
Set-TextValue $ws 'E7' '  -5.05%  '
Set-TextValue $ws 'E8' '  +0.07%  '
Set-TextValue $ws 'E9' '  -5.21%  '
Set-TextValue $ws 'E10' '  -6.47%  '
Set-TextValue $ws 'D11' '51.79'
Set-TextValue $ws 'E11' '  -2.00%  '
Set-TextValue $ws 'E12' '  -2.89%  '
Set-TextValue $ws 'E13' '  +0.51%  '
Set-TextValue $ws 'E14' '  -3.75%  '
Set-TextValue $ws 'D15' '2.665.57'
Set-TextValue $ws 'E15' '  -0.08%  '
Set-TextValue $ws 'D16' '15.69'
Set-TextValue $ws 'E16' '  +3.84%  '
Set-TextValue $ws 'D17' '2.295.77'
Set-TextValue $ws 'E17' '  -1.03%  '
Set-TextValue $ws 'D18' '0.812'
Set-TextValue $ws 'E18' '  +0.18%  '
Set-TextValue $ws 'D19' '42.671.23'
Set-TextValue $ws 'E19' '  -1.47%  '
Set-TextValue $ws 'E20' '  -2.55%  '
Set-TextValue $ws 'E21' '  -5.22%  '
Set-TextValue $ws 'D22' '6.09'
Set-TextValue $ws 'E22' '  -1.63%  '
Set-TextValue $ws 'D23' '69.27'
Set-TextValue $ws 'E23' '  +1.65%  '
Set-TextValue $ws 'D24' '234.97'
Set-TextValue $ws 'E24' '  -3.31%  '
Set-TextValue $ws 'E25' '  -2.24%  '
Set-TextValue $ws 'D26' '2.52'
Set-TextValue $ws 'E26' '  -3.36%  '
Set-TextValue $ws 'D28' '25.24'
Set-TextValue $ws 'E28' '  +1.28%  '
Set-TextValue $ws 'D29' '2.28'
Set-TextValue $ws 'E29' '  -1.35%  '
Set-TextValue $ws 'E30' '  -6.55%  '
Set-TextValue $ws 'E31' '  -4.63%  '
Set-TextValue $ws 'D32' '162.67'
Set-TextValue $ws 'E32' '  -2.88%  '
Set-TextValue $ws 'D33' '1.00'
Set-TextValue $ws 'E33' '  +0.01%  '
Set-TextValue $ws 'D34' '5.04'
Set-TextValue $ws 'E34' '  -4.57%  '
Set-TextValue $ws 'E35' '  +3.83%  '
Set-TextValue $ws 'E36' '  -3.58%  '
Set-TextValue $ws 'D37' '0.0715'
Set-TextValue $ws 'E37' '  -3.80%  '
Set-TextValue $ws 'D38' '16.99'
Set-TextValue $ws 'E38' '  -7.79%  '
Set-TextValue $ws 'E39' '  -5.52%  '
Set-TextValue $ws 'D40' '1.80'
Set-TextValue $ws 'E40' '  -4.01%  '
Set-TextValue $ws 'D41' '0.100'
Set-TextValue $ws 'E41' '  -5.28%  '
Set-TextValue $ws 'E42' '  -3.65%  '
Set-TextValue $ws 'E43' '  -8.62%  '
Set-TextValue $ws 'D44' '1.993.30'
Set-TextValue $ws 'E44' '  +0.20%  '
Set-TextValue $ws 'D45' '18.75'
Set-TextValue $ws 'E45' '  -1.67%  '
Set-TextValue $ws 'E46' '  -4.38%  '
Set-TextValue $ws 'D47' '10.23'
Set-TextValue $ws 'E47' '  +2.22%  '
Set-TextValue $ws 'D48' '2.88'
Set-TextValue $ws 'E48' '  -6.15%  '
Set-TextValue $ws 'D49' '55.47'
Set-TextValue $ws 'E49' '  -0.63%  '
Set-TextValue $ws 'E50' '  -2.57%  '
Set-TextValue $ws 'D51' '2.532.11'
Set-TextValue $ws 'E51' '  -0.14%  '
